# "split in correlative letters"
#
# The run "OCTAVA: FACULTADES PARA SUSTITUIR EL PODER" was wrapped in a
# reviewer comment ("CORRELATIVO INCORRECTO") and highlighted yellow to
# flag it. Resolve the review: drop the yellow highlight that marked the
# run and remove the comment (which also clears the comment anchors /
# reference around the run).

$d = $word.ActiveDocument

# --- 1. Clear the yellow highlight on the flagged heading run ----------
$rng = $d.Content
$found = $rng.Find.Execute("OCTAVA: FACULTADES PARA SUSTITUIR EL PODER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 0   # wdNoHighlight
}

# --- 2. Remove the "CORRELATIVO INCORRECTO" comment ---------------------
# Deleting the comment also removes its <w:commentRangeStart/End> and
# <w:commentReference> markers from the run.
$comments = $d.Comments
for ($i = $comments.Count; $i -ge 1; $i--) {
    $comments.Item($i).Delete()
}
